$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Volume 32   Number  21" -> "...  22" (issue number)
$ws.Range("A8").Characters(21, 2).Text = "22"

# Header: report week date range 5/19/2025-5/25/2025 -> 5/26/2025-6/1/2025
$ws.Range("C9").Characters(27, 9).Text = "5/26/2025"
$ws.Range("C9").Characters(47, 9).Text = "6/1/2025"

# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = 33.333333333333
$ws.Range("N15").Value = -20

# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -22.222222222222
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 80
$ws.Range("J16").Value = 101
$ws.Range("K16").Value = -20.79207920792
$ws.Range("L16").Value = -21.56862745098
$ws.Range("M16").Value = -13.043478260869
$ws.Range("N16").Value = -75.975975975976

# Row 17
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 50
$ws.Range("H17").Value = -36
$ws.Range("I17").Value = 133
$ws.Range("J17").Value = 204
$ws.Range("K17").Value = -34.803921568627
$ws.Range("L17").Value = -14.193548387096
$ws.Range("M17").Value = 27.884615384615
$ws.Range("N17").Value = -47.430830039525

# Row 18
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 54
$ws.Range("J18").Value = 41
$ws.Range("K18").Value = 31.70731707317
$ws.Range("L18").Value = 12.5
$ws.Range("M18").Value = 17.391304347826
$ws.Range("N18").Value = -74.285714285714

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -38.461538461538
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -18.918918918918
$ws.Range("I19").Value = 162
$ws.Range("J19").Value = 214
$ws.Range("K19").Value = -24.29906542056
$ws.Range("L19").Value = -4.705882352941
$ws.Range("M19").Value = 84.090909090909
$ws.Range("N19").Value = 18.248175182481

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 7
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 34
$ws.Range("J20").Value = 51
$ws.Range("K20").Value = -33.333333333333
$ws.Range("L20").Value = 25.925925925925
$ws.Range("M20").Value = 54.545454545454
$ws.Range("N20").Value = -73.846153846153

# Row 21
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -31.70731707317
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 126
$ws.Range("H21").Value = -26.984126984127
$ws.Range("I21").Value = 471
$ws.Range("J21").Value = 622
$ws.Range("K21").Value = -24.276527331189
$ws.Range("L21").Value = -8.0078125
$ws.Range("M21").Value = 29.395604395604
$ws.Range("N21").Value = -56.589861751152

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 11
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 10
$ws.Range("M22").Value = -38.888888888888

# Row 23
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 15.384615384615
$ws.Range("I23").Value = 71
$ws.Range("J23").Value = 74
$ws.Range("K23").Value = -4.054054054054
$ws.Range("L23").Value = 1.428571428571
$ws.Range("M23").Value = 77.5

# Row 24
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -37.037037037037
$ws.Range("F24").Value = 67
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = -22.988505747126
$ws.Range("I24").Value = 334
$ws.Range("J24").Value = 400
$ws.Range("K24").Value = -16.5
$ws.Range("L24").Value = -37.099811676082
$ws.Range("M24").Value = -10.933333333333

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = -37.5
$ws.Range("I25").Value = 86
$ws.Range("J25").Value = 95
$ws.Range("K25").Value = -9.473684210526
$ws.Range("L25").Value = -65.040650406504

# Row 26
$ws.Range("D26").Value = 19
$ws.Range("E26").Value = -42.105263157894
$ws.Range("G26").Value = 60
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 260
$ws.Range("J26").Value = 303
$ws.Range("K26").Value = -14.191419141914
$ws.Range("L26").Value = 25
$ws.Range("M26").Value = 40.54054054054

# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("L27").Value = -18.181818181818

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 24
$ws.Range("J28").Value = 24
$ws.Range("L28").Value = 4.347826086956

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("L29").Value = -76.470588235294
$ws.Range("M29").Value = -78.947368421052
$ws.Range("N29").Value = -77.777777777777

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("L30").Value = -75
$ws.Range("M30").Value = -75
$ws.Range("N30").Value = -77.777777777777

# Fix styles for cells that changed between text/number type
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("J14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("J14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null  # xlPasteFormats